# Re-position the "Data Objects" (left) and "Pipeline Objects" (right)
# boxes of the class-hierarchy diagram, along with their labels and the
# connector arrows attached to them -- matching the diagram layout shift
# recorded in the target OOXML.
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are single
# precision (32-bit) point values, and the host truncates point*12700
# down to whole EMUs. The literals below are the nearest points whose
# f32 rounding lands exactly on the target EMU offsets/extents, so the
# saved OOXML reproduces the diff precisely instead of drifting by an
# EMU here and there.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Dataset (id 5)
$sh = $s.Shapes.Item("TextBox 4")
$sh.Left = 177.8251190185547
$sh.Top = 121.33126068115234

# Subject (id 6)
$sh = $s.Shapes.Item("TextBox 5")
$sh.Left = 177.8251190185547
$sh.Top = 179.23126220703125

# Visit (id 7)
$sh = $s.Shapes.Item("TextBox 6")
$sh.Left = 177.8251190185547
$sh.Top = 237.1312713623047

# Trial (id 8)
$sh = $s.Shapes.Item("TextBox 7")
$sh.Left = 177.8251190185547
$sh.Top = 295.0312805175781

# Phase (id 9)
$sh = $s.Shapes.Item("TextBox 8")
$sh.Left = 177.8251190185547
$sh.Top = 352.9312744140625

# Project (id 10)
$sh = $s.Shapes.Item("TextBox 9")
$sh.Left = 401.1000061035156
$sh.Top = 155.70001220703125

# Analysis (id 11)
$sh = $s.Shapes.Item("TextBox 10")
$sh.Left = 401.1000061035156
$sh.Top = 213.60000610351562

# Process (id 12)
$sh = $s.Shapes.Item("TextBox 11")
$sh.Left = 401.1000061035156
$sh.Top = 271.5

# Dataset -> Subject connector (id 16)
$sh = $s.Shapes.Item("Straight Arrow Connector 15")
$sh.Left = 232.27513122558594
$sh.Top = 150.4125213623047

# Subject -> Visit connector (id 17)
$sh = $s.Shapes.Item("Straight Arrow Connector 16")
$sh.Left = 232.27513122558594
$sh.Top = 208.31253051757812

# Visit -> Trial connector (id 18)
$sh = $s.Shapes.Item("Straight Arrow Connector 17")
$sh.Left = 232.27513122558594
$sh.Top = 266.2125244140625

# Trial -> Phase connector (id 19)
$sh = $s.Shapes.Item("Straight Arrow Connector 18")
$sh.Left = 233.3251190185547
$sh.Top = 324.112548828125

# Project -> Analysis connector (id 20)
$sh = $s.Shapes.Item("Straight Arrow Connector 19")
$sh.Left = 454.6500244140625
$sh.Top = 184.78126525878906

# Analysis -> Process connector (id 21)
$sh = $s.Shapes.Item("Straight Arrow Connector 20")
$sh.Left = 454.6500244140625
$sh.Top = 242.6812744140625

# User -> Rectangle 25 (Data Objects box) connector (id 23)
$sh = $s.Shapes.Item("Straight Arrow Connector 22")
$sh.Left = 232.366943359375
$sh.Top = 65.86874389648438
$sh.Width = 135.0303955078125
$sh.Height = 43.462520599365234

# User -> Rectangle 29 (Pipeline Objects box) connector (id 25)
$sh = $s.Shapes.Item("Straight Arrow Connector 24")
$sh.Left = 367.3973388671875
$sh.Top = 65.86874389648438
$sh.Width = 155.26087951660156
$sh.Height = 76.48126220703125

# Data Objects group box (id 26)
$sh = $s.Shapes.Item("Rectangle 25")
$sh.Left = 163.90875244140625
$sh.Top = 109.33126068115234

# 'Data Objects' label (id 27)
$sh = $s.Shapes.Item("TextBox 26")
$sh.Left = 182.54087829589844
$sh.Top = 78.9000015258789

# Pipeline Objects group box (id 30)
$sh = $s.Shapes.Item("Rectangle 29")
$sh.Left = 391.0163269042969
$sh.Top = 142.35000610351562

# 'Pipeline Objects' label (id 31)
$sh = $s.Shapes.Item("TextBox 30")
$sh.Left = 454.4615783691406
$sh.Top = 113.26874542236328

# Subset (id 32)
$sh = $s.Shapes.Item("TextBox 31")
$sh.Left = 534.4500122070312
$sh.Top = 155.70001220703125

# Plot (id 33)
$sh = $s.Shapes.Item("TextBox 32")
$sh.Left = 534.4500122070312
$sh.Top = 213.24371337890625

# Logsheet (id 34)
$sh = $s.Shapes.Item("TextBox 33")
$sh.Left = 534.4500122070312
$sh.Top = 271.621826171875

# Rectangle 25 -> Variable connector (id 39)
$sh = $s.Shapes.Item("Straight Arrow Connector 38")
$sh.Left = 232.366943359375
$sh.Top = 392.8312683105469
$sh.Width = 135.0303955078125
$sh.Height = 42.20622253417969

# Rectangle 29 -> Variable connector (id 45)
$sh = $s.Shapes.Item("Straight Arrow Connector 44")
$sh.Left = 367.3973388671875
$sh.Top = 317.4000244140625
$sh.Width = 155.26087951660156
$sh.Height = 117.63748168945312
